$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("Power") values were stored as Watts (Voltage[V] * Current[mA] / 1000).
# They need to be rescaled to milliwatts by multiplying each value by 1000
# (i.e. Power = Voltage[V] * Current[mA], no /1000 division).
$lastRow = 321

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value()
    if ($val -ne $null) {
        $cell.Value = $val * 1000
    }
}
